$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values for columns B:E
$ws.Range("B2").Value = 461.62263204731903
$ws.Range("C2").Value = 389.61793576387959
$ws.Range("D2").Value = 457.90149340365843
$ws.Range("E2").Value = 385.30257189742673

# Update row 3 values for columns B:E
$ws.Range("B3").Value = 463.12862074754423
$ws.Range("C3").Value = 388.3339396377682
$ws.Range("D3").Value = 458.30069618978428
$ws.Range("E3").Value = 395.4634356969558

# Update the sheet selection to match the new active range
$ws.Range("B1:E3").Select()
